$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Size-class key data updates -------------------------------------------------
# B2 used to be a formula "(C2-1)+1/6"; it is now a plain static value.
$ws.Range("B2").Value = 0.26

# The "lower bound" values for the 2/3/4-way splits moved from 0.5 to the
# row's class-relative offset.
$ws.Range("B7").Value = 1.5
$ws.Range("B8").Value = 1.5
$ws.Range("B11").Value = 2.5
$ws.Range("B12").Value = 2.5
$ws.Range("B15").Value = 3.5
$ws.Range("B16").Value = 3.5
$ws.Range("B18").Value = 4.1666670000000003

# --- New helper column (E) used by Solver's objective cell ----------------------
$ws.Range("E2:E18").NumberFormat = "0.0"

# A handful of rows (the ones representing a whole integer size class) also pick
# up a blank, center-formatted cell in column F.
foreach ($r in @(4, 8, 12, 16)) {
    $ws.Cells.Item($r, 6).HorizontalAlignment = -4108
}

# Widen column E to fit its new contents.
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666

# --- Re-introduce the Solver add-in parameters (hidden defined names) -----------
function Add-HiddenName($name, $refersTo) {
    $n = $ws.Names.Add($name, $refersTo)
    $n.Visible = $false
}

Add-HiddenName 'solver_adj' '=Sheet1!$B$2'
Add-HiddenName 'solver_cvg' '0.0001'
Add-HiddenName 'solver_drv' '1'
Add-HiddenName 'solver_eng' '1'
Add-HiddenName 'solver_est' '1'
Add-HiddenName 'solver_itr' '2147483647'
Add-HiddenName 'solver_mip' '2147483647'
Add-HiddenName 'solver_mni' '30'
Add-HiddenName 'solver_mrt' '0.075'
Add-HiddenName 'solver_msl' '2'
Add-HiddenName 'solver_neg' '1'
Add-HiddenName 'solver_nod' '2147483647'
Add-HiddenName 'solver_num' '0'
Add-HiddenName 'solver_nwt' '1'
Add-HiddenName 'solver_opt' '=Sheet1!$E$2'
Add-HiddenName 'solver_pre' '0.000001'
Add-HiddenName 'solver_rbv' '1'
Add-HiddenName 'solver_rlx' '2'
Add-HiddenName 'solver_rsd' '0'
Add-HiddenName 'solver_scl' '1'
Add-HiddenName 'solver_sho' '2'
Add-HiddenName 'solver_ssz' '100'
Add-HiddenName 'solver_tim' '2147483647'
Add-HiddenName 'solver_tol' '0.01'
Add-HiddenName 'solver_typ' '3'
Add-HiddenName 'solver_val' '13'
Add-HiddenName 'solver_ver' '3'

# --- Final selection state --------------------------------------------------------
$ws.Range("B18").Select()
